$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.68"
$ws.Range("E2").Value = "'5.83%"
$ws.Range("D3").Value = "'35.28"
$ws.Range("E3").Value = "'13.67%"
$ws.Range("E4").Value = "'4.56%"
$ws.Range("D5").Value = "'0.07766"
$ws.Range("E5").Value = "'5.83%"
$ws.Range("D6").Value = "'2.395"
$ws.Range("E6").Value = "'6.30%"
$ws.Range("D7").Value = "'8.022"
$ws.Range("E7").Value = "'3.77%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9231"
$ws.Range("E8").Value = "'1.96%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.09988"
$ws.Range("E9").Value = "'14.37%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1801"
$ws.Range("E10").Value = "'6.80%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08594"
$ws.Range("E11").Value = "'3.83%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03309"
$ws.Range("E12").Value = "'6.31%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09896"
$ws.Range("E13").Value = "'-0.35%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001500"
$ws.Range("E14").Value = "'0.29%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005786"
$ws.Range("E15").Value = "'-0.59%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.471"
$ws.Range("E16").Value = "'-0.55%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'3.940"
$ws.Range("E17").Value = "'5.43%"
$ws.Range("D18").Value = "'2.164"
$ws.Range("E18").Value = "'4.34%"
$ws.Range("D19").Value = "'0.3365"
$ws.Range("E19").Value = "'1.10%"
$ws.Range("E20").Value = "'0.60%"
$ws.Range("D21").Value = "'4.304"
$ws.Range("E21").Value = "'3.19%"
$ws.Range("D22").Value = "'0.2384"
$ws.Range("E22").Value = "'12.24%"
$ws.Range("D23").Value = "'0.04571"
$ws.Range("E23").Value = "'0.47%"
$ws.Range("D24").Value = "'0.001216"
$ws.Range("E24").Value = "'0.57%"
$ws.Range("D25").Value = "'0.004458"
$ws.Range("E25").Value = "'7.60%"
$ws.Range("E26").Value = "'-0.28%"
$ws.Range("D39").Value = "'0.01787"
$ws.Range("E39").Value = "'13.81%"
$ws.Range("D40").Value = "'0.04747"
$ws.Range("E40").Value = "'6.79%"
$ws.Range("D41").Value = "'0.007738"
$ws.Range("E41").Value = "'5.23%"
$ws.Range("D42").Value = "'0.1412"
$ws.Range("E42").Value = "'6.73%"
$ws.Range("D43").Value = "'0.007089"
$ws.Range("E43").Value = "'-26.07%"
$ws.Range("E44").Value = "'-1.87%"
$ws.Range("D45").Value = "'0.009521"
$ws.Range("E45").Value = "'12.95%"
$ws.Range("E46").Value = "'0.09%"
$ws.Range("E47").Value = "'-0.29%"
$ws.Range("E48").Value = "'29.95%"
$ws.Range("D49").Value = "'0.001998"
$ws.Range("E49").Value = "'-0.23%"
$ws.Range("D50").Value = "'0.00002097"
$ws.Range("E50").Value = "'-0.29%"
$ws.Range("D51").Value = "'0.0001997"
$ws.Range("E51").Value = "'-0.29%"
